$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold + border + centered alignment) from the previous
# date header cell (Z1) onto the new one (AA1), then set its value.
$ws.Range("Z1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)
$ws.Range("AA1").Value = "22-03-2020"

# New confirmed-case counts per province for 22-03-2020
$ws.Range("AA2").Value = 41
$ws.Range("AA3").Value = 57
$ws.Range("AA4").Value = 34
$ws.Range("AA5").Value = 421
$ws.Range("AA6").Value = 56
$ws.Range("AA7").Value = 437
$ws.Range("AA8").Value = 1358
$ws.Range("AA9").Value = 499
$ws.Range("AA10").Value = 187
$ws.Range("AA11").Value = 344
$ws.Range("AA12").Value = 48
$ws.Range("AA13").Value = 567
